$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q16").Value = 509498
$ws.Range("R16").Value = 6551086
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()

$ws.Range("Q17").Value = 509365
$ws.Range("R17").Value = 6551082
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()

$ws.Range("Q18").Value = 509498
$ws.Range("R18").Value = 6551086
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
